$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns with
# the latest quote snapshot. The data cells (D2:E51) store their values as
# literal text (e.g. "1.001", "  -0.86%  ") rather than numbers/percentages,
# matching the feed that originally populated the sheet. Forcing the
# NumberFormat to text ("@") before assigning each Value keeps Excel from
# auto-converting numeric-looking strings into real numbers; ClearFormats()
# afterwards restores the cells to their original (unstyled) formatting.

$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.584.08'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '1.875.79'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '247.69'
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4742'
$ws.Range("E7").Value = '  -0.95%  '
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").Value = '0.06480'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").Value = '21.95'
$ws.Range("E10").Value = '  +2.52%  '
$ws.Range("D12").Value = '0.7390'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '1.877.20'
$ws.Range("E13").Value = '  -2.26%  '
$ws.Range("D14").Value = '95.80'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = '5.176'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").Value = '273.95'
$ws.Range("E16").Value = '  -2.85%  '
$ws.Range("D17").Value = '30.573.92'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("E18").Value = '  -2.76%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '0.000007471'
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").Value = '2.121.90'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '5.207'
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("D24").Value = '6.166'
$ws.Range("E24").Value = '  -1.39%  '
$ws.Range("D25").Value = '165.75'
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").Value = '9.187'
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").Value = '1.904'
$ws.Range("E28").Value = '  -4.49%  '
$ws.Range("D29").Value = '0.09890'
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("E30").Value = '  -2.64%  '
$ws.Range("D31").Value = '1.507'
$ws.Range("D32").Value = '4.244'
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("D33").Value = '4.087'
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").Value = '0.04770'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("D35").Value = '1.121'
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("D36").Value = '0.6940'
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("D37").Value = '2.720'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '0.01849'
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("D39").Value = '2.758'
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = '6.232'
$ws.Range("E40").Value = '  -3.54%  '
$ws.Range("D41").Value = '73.27'
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '0.4163'
$ws.Range("E44").Value = '  -1.59%  '
$ws.Range("D45").Value = '0.8345'
$ws.Range("E45").Value = '  -1.72%  '
$ws.Range("D46").Value = '101.50'
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("D47").Value = '9.339'
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("D48").Value = '35.33'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").Value = '6.964'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").Value = '922.87'
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").Value = '0.05668'
$ws.Range("E51").Value = '  +0.62%  '

$dataRange.ClearFormats()
